$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (B4 / row 8) entirely
$ws.Rows.Item(8).Delete()

# Row 2 (A1)
$ws.Range("B2").Value2 = 3.714285714285714
$ws.Range("C2").Value2 = 4.428571428571429
$ws.Range("D2").Value2 = 4
$ws.Range("E2").Value2 = 3.714285714285714
$ws.Range("F2").Value2 = 79.29000000000001
$ws.Range("G2").Value2 = "Should’ve done more about definitions`nThe topic is quite difficult The leader do his best`nIt was sure a very difficult topics topics discuss about. Great job though 👏`n"

# Row 3 (A2)
$ws.Range("B3").Value2 = 4.428571428571429
$ws.Range("C3").Value2 = 4.428571428571429
$ws.Range("D3").Value2 = 4.428571428571429
$ws.Range("E3").Value2 = 4.571428571428571
$ws.Range("F3").Value2 = 89.29000000000001
$ws.Range("G3").Value2 = "has a topic that really interesting to discuss`nThe topic is interested and the leader made it even interesting and easy for us to understand`nThe presentation slides are doing great`n"

# Row 4 (A3)
$ws.Range("B4").Value2 = 4.333333333333333
$ws.Range("C4").Value2 = 4.333333333333333
$ws.Range("D4").Value2 = 4.333333333333333
$ws.Range("E4").Value2 = 4.333333333333333
$ws.Range("F4").Value2 = 86.67
$ws.Range("G4").Value2 = "Good`nWe got a great talk during discussion`nDid well on explaining the topic.`ngood`nHe provide the information so clearly and have some little discusion make this session is not feel boring.`n"

# Row 5 (B1)
$ws.Range("B5").Value2 = 3.571428571428572
$ws.Range("C5").Value2 = 3.714285714285714
$ws.Range("D5").Value2 = 3.857142857142857
$ws.Range("E5").Value2 = 4
$ws.Range("F5").Value2 = 75.70999999999999
$ws.Range("G5").Value2 = "voice a bit unclear`nNice`nThe PowerPoint is clear and beautiful`n"

# Row 6 (B2)
$ws.Range("B6").Value2 = 2.714285714285714
$ws.Range("C6").Value2 = 2.714285714285714
$ws.Range("D6").Value2 = 2.857142857142857
$ws.Range("E6").Value2 = 2.857142857142857
$ws.Range("F6").Value2 = 55.71
$ws.Range("G6").Value2 = "I can see you are nervous and lack of practice, maybe next time you can put more effort in it and do it better.`nKeep up!`nactually has a good topic`nIt's a funny presentation. Love it.`n"

# Row 7 (B3)
$ws.Range("B7").Value2 = 4.166666666666667
$ws.Range("C7").Value2 = 4
$ws.Range("D7").Value2 = 4.333333333333333
$ws.Range("E7").Value2 = 4.333333333333333
$ws.Range("F7").Value2 = 84.17
$ws.Range("G7").Value2 = "Good job for leading the team, Information that provided so helpfull`nCould make sure that everyone has properly engaged, otherwise may finish early.`n"

